$d = $word.ActiveDocument
$W = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$apos = [char]0x2019

# ---------------------------------------------------------------------------
# Edit 1: paragraph 3 ("What to play is a software designed...") - split the
# "It's" (second occurrence, the contraction "It's purpose") into its own run
# wrapped in proofErr spellStart/spellEnd marks, matching Word's automatic
# spell-check annotation behaviour.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$full3 = $p3.Range

$xml3 = "<w:p $W>" + `
  "<w:r><w:t xml:space='preserve'>What to play is a software designed for gamers to better help them find new games. It offers tools such as price listings, sort by genre, and even a fully customizable profile with its own unique favorites list. </w:t></w:r>" + `
  "<w:proofErr w:type='spellStart'/>" + `
  "<w:r><w:t>It${apos}s</w:t></w:r>" + `
  "<w:proofErr w:type='spellEnd'/>" + `
  "<w:r><w:t xml:space='preserve'> purpose is to help you narrow down what game you want to play next, along with games you${apos}ve played and loved. </w:t></w:r>" + `
  "</w:p>"

$full3.InsertXML($xml3) | Out-Null

# ---------------------------------------------------------------------------
# Edit 2: the final paragraph ("Alright now you're in, but what can you do
# once you've logged in?") is rewritten into several runs, and a new bullet
# list describing the app's features is appended after it.
# ---------------------------------------------------------------------------
$p11 = $d.Paragraphs(11)
$full11 = $p11.Range

$xml11 = "<w:p $W>" + `
  "<w:r><w:t xml:space='preserve'>Alright now </w:t></w:r>" + `
  "<w:r><w:t>you are</w:t></w:r>" + `
  "<w:r><w:t xml:space='preserve'> in, but what can you do once </w:t></w:r>" + `
  "<w:r><w:t>you have</w:t></w:r>" + `
  "<w:r><w:t xml:space='preserve'> logged in?</w:t></w:r>" + `
  "<w:r><w:t xml:space='preserve'> Well, </w:t></w:r>" + `
  "<w:r><w:t>here</w:t></w:r>" + `
  "<w:r><w:t xml:space='preserve'> is</w:t></w:r>" + `
  "<w:r><w:t xml:space='preserve'> a list</w:t></w:r>" + `
  "<w:r><w:t>.</w:t></w:r>" + `
  "</w:p>"

$full11.InsertXML($xml11) | Out-Null

# Re-fetch the (now rewritten) last paragraph and append the new paragraphs
# describing each feature right after it.
$p11b = $d.Paragraphs(11)
$endPos = $p11b.Range.End - 1
$insertPoint = $d.Range($endPos, $endPos)

$newParasXml = `
  "<w:p $W><w:r><w:t>- Profile Customization</w:t></w:r></w:p>" + `
  "<w:p $W><w:r><w:t xml:space='preserve'>  - Our plain jane feature, it allows you to change username, password, </w:t></w:r><w:r><w:t>and all the other bells and whistles of your account.</w:t></w:r></w:p>" + `
  "<w:p $W><w:r><w:t>- Search Bar</w:t></w:r></w:p>" + `
  "<w:p $W><w:pPr><w:ind w:firstLine='105'/></w:pPr><w:r><w:t>- Our main feature of What to Play?</w:t></w:r><w:r><w:t>. It${apos}s a high functioning search bar with over 13000 games to search from!</w:t></w:r><w:r><w:t xml:space='preserve'> Each with their own description and images</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>" + `
  "<w:p $W><w:r><w:t>- Favorites List</w:t></w:r></w:p>" + `
  "<w:p $W><w:pPr><w:ind w:firstLine='105'/></w:pPr><w:r><w:t>- When searching our massive library of games, you are going to have some favorites. Why no pop those games into a condensed list of what you have really enjoyed playing. When look at a game, hit the button in the top right marked with a star and add it to your list</w:t></w:r></w:p>" + `
  "<w:p $W><w:r><w:t xml:space='preserve'>- </w:t></w:r><w:r><w:t>Wheel Spin</w:t></w:r></w:p>" + `
  "<w:p $W><w:r><w:t xml:space='preserve'>  - Lastly, our most fun feature. The wheel spin is a fun little wheel that you can </w:t></w:r><w:r><w:t>spin to be taken to a random game. With such an extensive library, you might not know what to play, instead let the system choose!</w:t></w:r></w:p>"

$insertPoint.InsertXML($newParasXml) | Out-Null

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
